$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Output $ws.Name
Write-Output $ws.Cells.Item(1,1).Value()
